$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 227, shifting all
# subsequent data down by two rows (227-278 -> 229-280).
$ws.Rows("227:228").Insert()

# New row 227 data
$ws.Cells.Item(227, 1).Value = 7
$ws.Cells.Item(227, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(227, 3).Value = "Ñuble"
$ws.Cells.Item(227, 4).Value = 44694
$ws.Cells.Item(227, 5).Value = 16
$ws.Cells.Item(227, 6).Value = 100112002
$ws.Cells.Item(227, 7).Value = "Pimiento"
$ws.Cells.Item(227, 8).Value = "Zafiro rojo"
$ws.Cells.Item(227, 9).Value = "Primera"
$ws.Cells.Item(227, 10).Value = 100
$ws.Cells.Item(227, 11).Value = 25000
$ws.Cells.Item(227, 12).Value = 26000
$ws.Cells.Item(227, 13).Value = 25500
$ws.Cells.Item(227, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(227, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(227, 16).Value = 1700
$ws.Cells.Item(227, 17).Value = 15
$ws.Cells.Item(227, 18).Value = "Hortaliza"

# New row 228 data
$ws.Cells.Item(228, 1).Value = 7
$ws.Cells.Item(228, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(228, 3).Value = "Ñuble"
$ws.Cells.Item(228, 4).Value = 44694
$ws.Cells.Item(228, 5).Value = 16
$ws.Cells.Item(228, 6).Value = 100112002
$ws.Cells.Item(228, 7).Value = "Pimiento"
$ws.Cells.Item(228, 8).Value = "Zafiro verde"
$ws.Cells.Item(228, 9).Value = "Primera"
$ws.Cells.Item(228, 10).Value = 100
$ws.Cells.Item(228, 11).Value = 15000
$ws.Cells.Item(228, 12).Value = 16000
$ws.Cells.Item(228, 13).Value = 15500
$ws.Cells.Item(228, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(228, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(228, 16).Value = 1033
$ws.Cells.Item(228, 17).Value = 15
$ws.Cells.Item(228, 18).Value = "Hortaliza"
